$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.040415947832038
$ws.Range("D2").Value = 1.048165430170753
$ws.Range("E2").Value = 1.055424580740648
$ws.Range("F2").Value = 1.06174510265412
$ws.Range("I2").Value = 1.04267959705488
$ws.Range("J2").Value = 1.045502674048733
$ws.Range("K2").Value = 1.050926036191133
$ws.Range("L2").Value = 1.058165075871118
$ws.Range("M2").Value = 1.064468329782145
$ws.Range("N2").Value = 1.046987407766715

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.041217604776512
$ws.Range("D3").Value = 1.048788668587357
$ws.Range("E3").Value = 1.056258035903114
$ws.Range("F3").Value = 1.062535270487316
$ws.Range("I3").Value = 1.042861480100807
$ws.Range("J3").Value = 1.045950677701138
$ws.Range("K3").Value = 1.051361754033107
$ws.Range("L3").Value = 1.058811918596419
$ws.Range("M3").Value = 1.065073240252192
$ws.Range("N3").Value = 1.047436047635693

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.041736979809363
$ws.Range("D4").Value = 1.04919255604656
$ws.Range("E4").Value = 1.056798761180355
$ws.Range("F4").Value = 1.063047670460327
$ws.Range("I4").Value = 1.042978357773788
$ws.Range("J4").Value = 1.04624052594192
$ws.Range("K4").Value = 1.051643609634729
$ws.Range("L4").Value = 1.059231241356021
$ws.Range("M4").Value = 1.065465105241255
$ws.Range("N4").Value = 1.047726307494223

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.041955478738236
$ws.Range("D5").Value = 1.049362494638217
$ws.Range("E5").Value = 1.057026420937879
$ws.Range("F5").Value = 1.063263346730049
$ws.Range("I5").Value = 1.043027297806489
$ws.Range("J5").Value = 1.04636236723707
$ws.Range("K5").Value = 1.051762080515293
$ws.Range("L5").Value = 1.059407707811747
$ws.Range("M5").Value = 1.065629950718166
$ws.Range("N5").Value = 1.047848321817987

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.041992174628584
$ws.Range("D6").Value = 1.049391036483504
$ws.Range("E6").Value = 1.057064665823094
$ws.Range("F6").Value = 1.063299575105909
$ws.Range("I6").Value = 1.043035503576658
$ws.Range("J6").Value = 1.046382824249078
$ws.Range("K6").Value = 1.051781971014188
$ws.Range("L6").Value = 1.059437347997364
$ws.Range("M6").Value = 1.065657635132789
$ws.Range("N6").Value = 1.047868807881299

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.041739898801196
$ws.Range("D7").Value = 1.049194826209662
$ws.Range("E7").Value = 1.056801801852327
$ws.Range("F7").Value = 1.063050551305126
$ws.Range("I7").Value = 1.042979012481392
$ws.Range("J7").Value = 1.046242154035762
$ws.Range("K7").Value = 1.051645192733485
$ws.Range("L7").Value = 1.059233598591868
$ws.Range("M7").Value = 1.065467307502415
$ws.Range("N7").Value = 1.047727937900145

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.040686736008059
$ws.Range("D8").Value = 1.048375928981011
$ws.Range("E8").Value = 1.0557059547098
$ws.Range("F8").Value = 1.062011912797706
$ws.Range("I8").Value = 1.042741233129931
$ws.Range("J8").Value = 1.045654086491492
$ws.Range("K8").Value = 1.05107330544564
$ws.Range("L8").Value = 1.058383518465967
$ws.Range("M8").Value = 1.064672668441547
$ws.Range("N8").Value = 1.047139035232511

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.038835981786538
$ws.Range("D9").Value = 1.046937685510188
$ws.Range("E9").Value = 1.053785930256196
$ws.Range("F9").Value = 1.06019028068888
$ws.Range("I9").Value = 1.04231604356288
$ws.Range("J9").Value = 1.044617587357594
$ws.Range("K9").Value = 1.050064991171079
$ws.Range("L9").Value = 1.056891559140427
$ws.Range("M9").Value = 1.063275914679038
$ws.Range("N9").Value = 1.046101064150954

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.037605650610359
$ws.Range("D10").Value = 1.045982166684803
$ws.Range("E10").Value = 1.052513432535449
$ws.Range("F10").Value = 1.0589817528378
$ws.Range("I10").Value = 1.042028465684319
$ws.Range("J10").Value = 1.043926498849202
$ws.Range("K10").Value = 1.049392475213435
$ws.Range("L10").Value = 1.05590104983795
$ws.Range("M10").Value = 1.062347200971661
$ws.Range("N10").Value = 1.045408994217611

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.03707375765913
$ws.Range("D11").Value = 1.045569225157236
$ws.Range("E11").Value = 1.051964236581255
$ws.Range("F11").Value = 1.058459871042979
$ws.Range("I11").Value = 1.041902974020509
$ws.Range("J11").Value = 1.043627245329493
$ws.Range("K11").Value = 1.049101212776625
$ws.Range("L11").Value = 1.055473149763529
$ws.Range("M11").Value = 1.061945661499806
$ws.Range("N11").Value = 1.04510931572358

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.036876318254227
$ws.Range("D12").Value = 1.045415963060551
$ws.Range("E12").Value = 1.051760513922685
$ws.Range("F12").Value = 1.058266236195366
$ws.Range("I12").Value = 1.041856216072839
$ws.Range("J12").Value = 1.043516089437607
$ws.Range("K12").Value = 1.048993017516726
$ws.Range("L12").Value = 1.05531436018036
$ws.Range("M12").Value = 1.06179660369647
$ws.Range("N12").Value = 1.044998001977576

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.03691866379537
$ws.Range("D13").Value = 1.045448832714077
$ws.Range("E13").Value = 1.051804200728426
$ws.Range("F13").Value = 1.058307761750944
$ws.Range("I13").Value = 1.041866252349694
$ws.Range("J13").Value = 1.043539932723702
$ws.Range("K13").Value = 1.049016226080977
$ws.Range("L13").Value = 1.055348414198293
$ws.Range("M13").Value = 1.06182857291764
$ws.Range("N13").Value = 1.045021879123873

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.037057434607641
$ws.Range("D14").Value = 1.04555655394891
$ws.Range("E14").Value = 1.051947391217706
$ws.Range("F14").Value = 1.058443860721903
$ws.Range("I14").Value = 1.041899111942568
$ws.Range("J14").Value = 1.043618057137196
$ws.Range("K14").Value = 1.049092269466124
$ws.Range("L14").Value = 1.055460021065101
$ws.Range("M14").Value = 1.061933338451824
$ws.Range("N14").Value = 1.045100114482995

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.037142953076064
$ws.Range("D15").Value = 1.045622940882818
$ws.Range("E15").Value = 1.052035651866177
$ws.Range("F15").Value = 1.058527744381403
$ws.Range("I15").Value = 1.041919338663305
$ws.Range("J15").Value = 1.043666192231891
$ws.Range("K15").Value = 1.049139121363197
$ws.Range("L15").Value = 1.055528805874917
$ws.Range("M15").Value = 1.061997900176159
$ws.Range("N15").Value = 1.045148317935046

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.037640968628441
$ws.Range("D16").Value = 1.046009589366045
$ws.Range("E16").Value = 1.052549919052821
$ws.Range("F16").Value = 1.059016418487495
$ws.Range("I16").Value = 1.042036773811172
$ws.Range("J16").Value = 1.043946359272968
$ws.Range("K16").Value = 1.049411804225503
$ws.Range("L16").Value = 1.055929469313179
$ws.Range("M16").Value = 1.062373862607157
$ws.Range("N16").Value = 1.045428882845457

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.037953589105005
$ws.Range("D17").Value = 1.046252340466694
$ws.Range("E17").Value = 1.052872989758828
$ws.Range("F17").Value = 1.059323332190221
$ws.Range("I17").Value = 1.042110179003479
$ws.Range("J17").Value = 1.044122099555307
$ws.Range("K17").Value = 1.049582836269011
$ws.Range("L17").Value = 1.05618106307485
$ws.Range("M17").Value = 1.062609855688912
$ws.Range("N17").Value = 1.045604872699154

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.038136017091262
$ws.Range("D18").Value = 1.046394010564051
$ws.Range("E18").Value = 1.053061605347308
$ws.Range("F18").Value = 1.059502486504918
$ws.Range("I18").Value = 1.04215290150887
$ws.Range("J18").Value = 1.044224605062374
$ws.Range("K18").Value = 1.049682590620803
$ws.Range("L18").Value = 1.056327909535758
$ws.Range("M18").Value = 1.062747564205833
$ws.Range("N18").Value = 1.045707523775799

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.038198234130971
$ws.Range("D19").Value = 1.046442329491028
$ws.Range("E19").Value = 1.053125947821792
$ws.Range("F19").Value = 1.059563596625641
$ws.Range("I19").Value = 1.042167452896714
$ws.Range("J19").Value = 1.044259556610482
$ws.Range("K19").Value = 1.049716603260352
$ws.Range("L19").Value = 1.056377996626434
$ws.Range("M19").Value = 1.062794528993075
$ws.Range("N19").Value = 1.045742524959115

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.037920039417341
$ws.Range("D20").Value = 1.046226287542266
$ws.Range("E20").Value = 1.05283830932561
$ws.Range("F20").Value = 1.059290389094985
$ws.Range("I20").Value = 1.042102312986878
$ws.Range("J20").Value = 1.044103244362035
$ws.Range("K20").Value = 1.049564486733715
$ws.Range("L20").Value = 1.056154059494893
$ws.Range("M20").Value = 1.062584529879287
$ws.Range("N20").Value = 1.045585990729346

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.037016566477614
$ws.Range("D21").Value = 1.045524829325622
$ws.Range("E21").Value = 1.051905217629418
$ws.Range("F21").Value = 1.05840377700697
$ws.Range("I21").Value = 1.041889439601364
$ws.Range("J21").Value = 1.043595051428652
$ws.Range("K21").Value = 1.049069876778343
$ws.Range("L21").Value = 1.055427151426631
$ws.Range("M21").Value = 1.061902485064164
$ws.Range("N21").Value = 1.045077076103706

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.036449266294617
$ws.Range("D22").Value = 1.045084505976812
$ws.Range("E22").Value = 1.051320127396067
$ws.Range("F22").Value = 1.057847575444077
$ws.Range("I22").Value = 1.041754760516084
$ws.Range("J22").Value = 1.043275532120103
$ws.Range("K22").Value = 1.04875885360766
$ws.Range("L22").Value = 1.054970993209793
$ws.Range("M22").Value = 1.061474188608213
$ws.Range("N22").Value = 1.04475710304109

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.036749931235537
$ws.Range("D23").Value = 1.0453178616555
$ws.Range("E23").Value = 1.051630144169312
$ws.Range("F23").Value = 1.058142309534849
$ws.Range("I23").Value = 1.041826235545336
$ws.Range("J23").Value = 1.043444914737952
$ws.Range("K23").Value = 1.048923736420249
$ws.Range("L23").Value = 1.05521272752099
$ws.Range("M23").Value = 1.061701185637372
$ws.Range("N23").Value = 1.044926726201684

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.037935198824703
$ws.Range("D24").Value = 1.046238059498214
$ws.Range("E24").Value = 1.05285397938351
$ws.Range("F24").Value = 1.059305274239472
$ws.Range("I24").Value = 1.042105867590193
$ws.Range("J24").Value = 1.044111764214989
$ws.Range("K24").Value = 1.049572778116802
$ws.Range("L24").Value = 1.056166260953696
$ws.Range("M24").Value = 1.062595973344327
$ws.Range("N24").Value = 1.045594522681468

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.039313836221037
$ws.Range("D25").Value = 1.047308930610489
$ws.Range("E25").Value = 1.054280986476997
$ws.Range("F25").Value = 1.060660186532083
$ws.Range("I25").Value = 1.042426694412175
$ws.Range("J25").Value = 1.044885568666081
$ws.Range("K25").Value = 1.050325724501479
$ws.Range("L25").Value = 1.057276545915833
$ws.Range("M25").Value = 1.063636583613716
$ws.Range("N25").Value = 1.046369426023636
